# [feature/ResolveReport] mail template fixes
# Duplicate row 3 ("Appreciations" sheet) into a new row 9, matching the
# target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Appreciations")

$ws.Range("A9").Value = "corevalue 1 updated"
$ws.Range("B9").Value = "new desc 1"
$ws.Range("C9").Value = "Great teamwork dskafkdfgkjahda a dhsfkjadhkfha hdkjfhakdhfkg sadhfkahdsjkfhakj aksjdhfkjahdkjfhakjsd akdshfkjahdfk"
$ws.Range("D9").Value = "Samnit"
$ws.Range("E9").Value = "Patil"
$ws.Range("F9").Value = "Software Engineer"
$ws.Range("G9").Value = "Sharyu"
$ws.Range("H9").Value = "Marwadi"
$ws.Range("I9").Value = "Trainee"
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
